$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2502").Value = '2025-07-25T15:22:09+00:00'
$ws.Range("B2502").Value = 'EXTERNAL:- Week 33 - BBC Radio Ulster/BBC Radio Foyle Billings Schedule - Press Issue'
$ws.Range("C2502").Value = 'pressportal@bbc.co.uk'
$ws.Range("D2502").Value = $true

$ws.Range("A2503").Value = '2025-07-25T15:22:07+00:00'
$ws.Range("B2503").Value = 'EXTERNAL:- Week 33 - BBC Radio Ulster/BBC Radio Foyle Billings Schedule - Press Issue'
$ws.Range("C2503").Value = 'pressportal@bbc.co.uk'
$ws.Range("D2503").Value = $true

$ws.Range("A2504").Value = '2025-07-25T15:22:06+00:00'
$ws.Range("B2504").Value = 'EXTERNAL:- Weekly Schedule released 25.07.25'
$ws.Range("C2504").Value = 'Bethan.Sloman@qvc.com'
$ws.Range("D2504").Value = $true

$ws.Range("A2505").Value = '2025-07-25T15:22:02+00:00'
$ws.Range("B2505").Value = 'EXTERNAL:- Weekly Schedule for w/c 21.07.25 released 25.07.25'
$ws.Range("C2505").Value = 'Bethan.Sloman@qvc.com'
$ws.Range("D2505").Value = $true

$ws.Range("A2506").Value = '2025-07-25T15:21:56+00:00'
$ws.Range("B2506").Value = 'EXTERNAL:- Week 33 - BBC Radio Ulster/BBC Radio Foyle Billings Schedule - Press Issue'
$ws.Range("C2506").Value = 'elaine.fullerton@bbc.co.uk'
$ws.Range("D2506").Value = $true

$ws.Range("A2507").Value = '2025-07-25T15:21:17+00:00'
$ws.Range("B2507").Value = 'EXTERNAL:- ITV4 Post Press Change - Monday 28th July'
$ws.Range("C2507").Value = 'itv2-4-scheduling@itv.com'
$ws.Range("D2507").Value = $false

$ws.Range("A2508").Value = '2025-07-25T15:17:44+00:00'
$ws.Range("B2508").Value = 'EXTERNAL:- Radio 4 Forward Schedule Week 33 2025 Issue 1'
$ws.Range("C2508").Value = 'patricia.hetherington@bbc.co.uk'
$ws.Range("D2508").Value = $true

$ws.Range("A2509").Value = '2025-07-25T15:16:43+00:00'
$ws.Range("B2509").Value = 'EXTERNAL:- BBC Radio 5 Live - Wk30 - 2025-07-27 - Sunday'
$ws.Range("C2509").Value = 'pressportal@bbc.co.uk'
$ws.Range("D2509").Value = $true

$ws.Range("A2510").Value = '2025-07-25T15:16:40+00:00'
$ws.Range("B2510").Value = 'EXTERNAL:- BBC Asian Network - Wk30 - 2025-07-27 - Sunday'
$ws.Range("C2510").Value = 'pressportal@bbc.co.uk'
$ws.Range("D2510").Value = $true

$ws.Range("A2511").Value = '2025-07-25T15:16:37+00:00'
$ws.Range("B2511").Value = 'EXTERNAL:- BBC Radio 4 FM - Wk30 - 2025-07-26 - Saturday'
$ws.Range("C2511").Value = 'pressportal@bbc.co.uk'
$ws.Range("D2511").Value = $true

$ws.Range("A2512").Value = '2025-07-25T15:16:37+00:00'
$ws.Range("B2512").Value = 'EXTERNAL:- BBC Radio Scotland Extra - Wk30 - 2025-07-27 - Sunday'
$ws.Range("C2512").Value = 'pressportal@bbc.co.uk'
$ws.Range("D2512").Value = $true

$ws.Range("A2513").Value = '2025-07-25T15:16:35+00:00'
$ws.Range("B2513").Value = 'EXTERNAL:- BBC Radio Orkney - Wk30 - 2025-07-26 - Saturday'
$ws.Range("C2513").Value = 'pressportal@bbc.co.uk'
$ws.Range("D2513").Value = $true

$ws.Range("A2514").Value = '2025-07-25T15:16:34+00:00'
$ws.Range("B2514").Value = 'EXTERNAL:- BBC Radio Orkney - Wk29 - 2025-07-25 - Friday'
$ws.Range("C2514").Value = 'pressportal@bbc.co.uk'
$ws.Range("D2514").Value = $true

$ws.Range("A2515").Value = '2025-07-25T15:16:33+00:00'
$ws.Range("B2515").Value = 'EXTERNAL:- BBC World Service UK Schedule - Wk29 - 2025-07-25 - Friday'
$ws.Range("C2515").Value = 'pressportal@bbc.co.uk'
$ws.Range("D2515").Value = $true

$ws.Range("A2516").Value = '2025-07-25T15:16:25+00:00'
$ws.Range("B2516").Value = 'EXTERNAL:- Pop Schedule for Week 33-35'
$ws.Range("C2516").Value = 'schedules@globallistings-distribution.com'
$ws.Range("D2516").Value = $true

$ws.Range("A2517").Value = '2025-07-25T15:16:24+00:00'
$ws.Range("B2517").Value = 'EXTERNAL:- BBC Radio Wales FM - Wk30 - 2025-08-01 - Friday'
$ws.Range("C2517").Value = 'pressportal@bbc.co.uk'
$ws.Range("D2517").Value = $true

$ws.Range("A2518").Value = '2025-07-25T15:16:19+00:00'
$ws.Range("B2518").Value = 'EXTERNAL:- BBC World Service UK Schedule - Wk30 - 2025-07-31 - Thursday'
$ws.Range("C2518").Value = 'pressportal@bbc.co.uk'
$ws.Range("D2518").Value = $true

$ws.Range("A2519").Value = '2025-07-25T15:16:16+00:00'
$ws.Range("B2519").Value = 'EXTERNAL:- BBC Asian Network - Wk30 - 2025-07-30 - Wednesday'
$ws.Range("C2519").Value = 'pressportal@bbc.co.uk'
$ws.Range("D2519").Value = $true

$ws.Range("A2520").Value = '2025-07-25T15:16:15+00:00'
$ws.Range("B2520").Value = 'EXTERNAL:- BBC Radio Scotland Extra - Wk30 - 2025-07-30 - Wednesday'
$ws.Range("C2520").Value = 'pressportal@bbc.co.uk'
$ws.Range("D2520").Value = $true

$ws.Range("A2521").Value = '2025-07-25T15:16:10+00:00'
$ws.Range("B2521").Value = 'EXTERNAL:- BBC World Service UK Schedule - Wk30 - 2025-07-27 - Sunday'
$ws.Range("C2521").Value = 'pressportal@bbc.co.uk'
$ws.Range("D2521").Value = $true

$ws.Range("A2522").Value = '2025-07-25T15:16:09+00:00'
$ws.Range("B2522").Value = 'EXTERNAL:- BBC Radio Wales FM - Wk30 - 2025-07-27 - Sunday'
$ws.Range("C2522").Value = 'pressportal@bbc.co.uk'
$ws.Range("D2522").Value = $true

$ws.Range("A2523").Value = '2025-07-25T15:16:07+00:00'
$ws.Range("B2523").Value = 'EXTERNAL:- BBC World Service UK Schedule - Wk30 - 2025-07-26 - Saturday'
$ws.Range("C2523").Value = 'pressportal@bbc.co.uk'
$ws.Range("D2523").Value = $true

$ws.Range("A2524").Value = '2025-07-25T15:16:06+00:00'
$ws.Range("B2524").Value = 'EXTERNAL:- BBC Radio Shetland - Wk30 - 2025-08-01 - Friday'
$ws.Range("C2524").Value = 'pressportal@bbc.co.uk'
$ws.Range("D2524").Value = $true

$ws.Range("A2525").Value = '2025-07-25T15:16:06+00:00'
$ws.Range("B2525").Value = 'EXTERNAL:- BBC Radio Scotland - Wk30 - 2025-08-01 - Friday'
$ws.Range("C2525").Value = 'pressportal@bbc.co.uk'
$ws.Range("D2525").Value = $true

$ws.Range("A2526").Value = '2025-07-25T15:16:06+00:00'
$ws.Range("B2526").Value = 'EXTERNAL:- BBC Radio Wales Extra - Wk30 - 2025-07-27 - Sunday'
$ws.Range("C2526").Value = 'pressportal@bbc.co.uk'
$ws.Range("D2526").Value = $true

$ws.Range("A2527").Value = '2025-07-25T15:16:05+00:00'
$ws.Range("B2527").Value = 'EXTERNAL:- BBC Radio Scotland Extra - Wk30 - 2025-08-01 - Friday'
$ws.Range("C2527").Value = 'pressportal@bbc.co.uk'
$ws.Range("D2527").Value = $true

$ws.Range("A2528").Value = '2025-07-25T15:16:04+00:00'
$ws.Range("B2528").Value = 'EXTERNAL:- BBC Radio Shetland - Wk31 - 2025-08-02 - Saturday'
$ws.Range("C2528").Value = 'pressportal@bbc.co.uk'
$ws.Range("D2528").Value = $true

$ws.Range("A2529").Value = '2025-07-25T15:16:03+00:00'
$ws.Range("B2529").Value = 'EXTERNAL:- BBC Radio Shetland - Wk30 - 2025-07-26 - Saturday'
$ws.Range("C2529").Value = 'pressportal@bbc.co.uk'
$ws.Range("D2529").Value = $true

$ws.Range("A2530").Value = '2025-07-25T15:16:03+00:00'
$ws.Range("B2530").Value = 'EXTERNAL:- BBC Radio 1 - Wk30 - 2025-07-26 - Saturday'
$ws.Range("C2530").Value = 'pressportal@bbc.co.uk'
$ws.Range("D2530").Value = $true

$ws.Range("A2531").Value = '2025-07-25T15:16:01+00:00'
$ws.Range("B2531").Value = 'EXTERNAL:- BBC Radio 5 Live - Wk30 - 2025-07-26 - Saturday'
$ws.Range("C2531").Value = 'pressportal@bbc.co.uk'
$ws.Range("D2531").Value = $true

$tbl = $ws.ListObjects.Item(1)
$tbl.Resize($ws.Range("A1:D2531"))
